# Rebuild the living_rooms memory-block trial table (rows 2-42) to match the
# newly generated 20-version stimulus set: renumber trial_total (F) sequentially
# from 365, and reshuffle each row's condition/stimulus/rating columns (G..V)
# according to the new row order (the catch trial moves from row 21 to row 23).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 53
$ws.Range("B2").Value = 'memory'
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 365
$ws.Range("G2").Value = 'living_rooms'
$ws.Range("H2").Value = 'living_rooms'
$ws.Range("I2").Value = 'target'
$ws.Range("J2").Value = 'old'
$ws.Range("K2").Value = 'j'
$ws.Range("L2").Value = 'stimuli/img_wgddx.png'
$ws.Range("M2").Value = 45.6304347826087
$ws.Range("N2").Value = 34.30434782608695
$ws.Range("O2").Value = 39.96739130434783
$ws.Range("P2").Value = 46
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 3
$ws.Range("S2").Value = 3
$ws.Range("T2").Value = 3
$ws.Range("U2").Value = 3
$ws.Range("V2").Value = 4

# Row 3
$ws.Range("A3").Value = 53
$ws.Range("B3").Value = 'memory'
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 366
$ws.Range("G3").Value = 'living_rooms'
$ws.Range("H3").Value = 'living_rooms'
$ws.Range("I3").Value = 'target'
$ws.Range("J3").Value = 'old'
$ws.Range("K3").Value = 'j'
$ws.Range("L3").Value = 'stimuli/img_a9he3.png'
$ws.Range("M3").Value = 83.06521739130434
$ws.Range("N3").Value = 63.95652173913044
$ws.Range("O3").Value = 73.51086956521739
$ws.Range("P3").Value = 46
$ws.Range("Q3").Value = 8
$ws.Range("R3").Value = 8
$ws.Range("S3").Value = 8
$ws.Range("T3").Value = 8
$ws.Range("U3").Value = 8
$ws.Range("V3").Value = 8

# Row 4
$ws.Range("A4").Value = 53
$ws.Range("B4").Value = 'memory'
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 367
$ws.Range("G4").Value = 'living_rooms'
$ws.Range("H4").Value = 'living_rooms'
$ws.Range("I4").ClearContents()
$ws.Range("J4").Value = 'new'
$ws.Range("K4").Value = 'f'
$ws.Range("L4").Value = 'stimuli/img_j4ttn.png'
$ws.Range("M4").Value = 12.61904761904762
$ws.Range("N4").Value = 11.42857142857143
$ws.Range("O4").Value = 12.02380952380952
$ws.Range("P4").Value = 42
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = 1
$ws.Range("V4").Value = 1

# Row 5
$ws.Range("A5").Value = 53
$ws.Range("B5").Value = 'memory'
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 368
$ws.Range("G5").Value = 'living_rooms'
$ws.Range("H5").Value = 'living_rooms'
$ws.Range("I5").Value = 'target'
$ws.Range("J5").Value = 'old'
$ws.Range("K5").Value = 'j'
$ws.Range("L5").Value = 'stimuli/img_pbsj1.png'
$ws.Range("M5").Value = 73.88636363636364
$ws.Range("N5").Value = 51.52272727272727
$ws.Range("O5").Value = 62.70454545454545
$ws.Range("P5").Value = 44
$ws.Range("Q5").Value = 6
$ws.Range("R5").Value = 6
$ws.Range("S5").Value = 6
$ws.Range("T5").Value = 6
$ws.Range("U5").Value = 6
$ws.Range("V5").Value = 6

# Row 6
$ws.Range("A6").Value = 53
$ws.Range("B6").Value = 'memory'
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 369
$ws.Range("G6").Value = 'living_rooms'
$ws.Range("H6").Value = 'living_rooms'
$ws.Range("I6").ClearContents()
$ws.Range("J6").Value = 'new'
$ws.Range("K6").Value = 'f'
$ws.Range("L6").Value = 'stimuli/img_vgh2g.png'
$ws.Range("M6").Value = 93.81395348837209
$ws.Range("N6").Value = 78.27906976744185
$ws.Range("O6").Value = 86.04651162790697
$ws.Range("P6").Value = 43
$ws.Range("Q6").Value = 10
$ws.Range("R6").Value = 10
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = 10
$ws.Range("U6").Value = 10
$ws.Range("V6").Value = 10

# Row 7
$ws.Range("A7").Value = 53
$ws.Range("B7").Value = 'memory'
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 370
$ws.Range("G7").Value = 'living_rooms'
$ws.Range("H7").Value = 'living_rooms'
$ws.Range("I7").ClearContents()
$ws.Range("J7").Value = 'new'
$ws.Range("K7").Value = 'f'
$ws.Range("L7").Value = 'stimuli/img_j856a.png'
$ws.Range("M7").Value = 38.225
$ws.Range("N7").Value = 25.875
$ws.Range("O7").Value = 32.05
$ws.Range("P7").Value = 40
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 2
$ws.Range("S7").Value = 2
$ws.Range("T7").Value = 3
$ws.Range("U7").Value = 3
$ws.Range("V7").Value = 2

# Row 8
$ws.Range("A8").Value = 53
$ws.Range("B8").Value = 'memory'
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 371
$ws.Range("G8").Value = 'living_rooms'
$ws.Range("H8").Value = 'living_rooms'
$ws.Range("I8").ClearContents()
$ws.Range("J8").Value = 'new'
$ws.Range("K8").Value = 'f'
$ws.Range("L8").Value = 'stimuli/img_6a0hu.png'
$ws.Range("M8").Value = 61.275
$ws.Range("N8").Value = 42.025
$ws.Range("O8").Value = 51.65
$ws.Range("P8").Value = 40
$ws.Range("Q8").Value = 4
$ws.Range("R8").Value = 4
$ws.Range("S8").Value = 4
$ws.Range("T8").Value = 5
$ws.Range("U8").Value = 4
$ws.Range("V8").Value = 5

# Row 9
$ws.Range("A9").Value = 53
$ws.Range("B9").Value = 'memory'
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 372
$ws.Range("G9").Value = 'living_rooms'
$ws.Range("H9").Value = 'living_rooms'
$ws.Range("I9").ClearContents()
$ws.Range("J9").Value = 'new'
$ws.Range("K9").Value = 'f'
$ws.Range("L9").Value = 'stimuli/img_c89x3.png'
$ws.Range("M9").Value = 72.8695652173913
$ws.Range("N9").Value = 49.65217391304348
$ws.Range("O9").Value = 61.26086956521739
$ws.Range("P9").Value = 46
$ws.Range("Q9").Value = 6
$ws.Range("R9").Value = 6
$ws.Range("S9").Value = 6
$ws.Range("T9").Value = 6
$ws.Range("U9").Value = 6
$ws.Range("V9").Value = 5

# Row 10
$ws.Range("A10").Value = 53
$ws.Range("B10").Value = 'memory'
$ws.Range("C10").Value = 6
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 373
$ws.Range("G10").Value = 'living_rooms'
$ws.Range("H10").Value = 'living_rooms'
$ws.Range("I10").ClearContents()
$ws.Range("J10").Value = 'new'
$ws.Range("K10").Value = 'f'
$ws.Range("L10").Value = 'stimuli/img_5jp4f.png'
$ws.Range("M10").Value = 84.85714285714286
$ws.Range("N10").Value = 67.83333333333333
$ws.Range("O10").Value = 76.3452380952381
$ws.Range("P10").Value = 42
$ws.Range("Q10").Value = 9
$ws.Range("R10").Value = 9
$ws.Range("S10").Value = 9
$ws.Range("T10").Value = 8
$ws.Range("U10").Value = 8
$ws.Range("V10").Value = 9

# Row 11
$ws.Range("A11").Value = 53
$ws.Range("B11").Value = 'memory'
$ws.Range("C11").Value = 6
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 374
$ws.Range("G11").Value = 'living_rooms'
$ws.Range("H11").Value = 'living_rooms'
$ws.Range("I11").Value = 'target'
$ws.Range("J11").Value = 'old'
$ws.Range("K11").Value = 'j'
$ws.Range("L11").Value = 'stimuli/img_i6wsx.png'
$ws.Range("M11").Value = 79.07142857142857
$ws.Range("N11").Value = 58
$ws.Range("O11").Value = 68.53571428571428
$ws.Range("P11").Value = 42
$ws.Range("Q11").Value = 7
$ws.Range("R11").Value = 7
$ws.Range("S11").Value = 7
$ws.Range("T11").Value = 7
$ws.Range("U11").Value = 7
$ws.Range("V11").Value = 7

# Row 12
$ws.Range("A12").Value = 53
$ws.Range("B12").Value = 'memory'
$ws.Range("C12").Value = 6
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 11
$ws.Range("F12").Value = 375
$ws.Range("G12").Value = 'living_rooms'
$ws.Range("H12").Value = 'living_rooms'
$ws.Range("I12").ClearContents()
$ws.Range("J12").Value = 'new'
$ws.Range("K12").Value = 'f'
$ws.Range("L12").Value = 'stimuli/img_x4bln.png'
$ws.Range("M12").Value = 76.34042553191489
$ws.Range("N12").Value = 59.51063829787234
$ws.Range("O12").Value = 67.92553191489361
$ws.Range("P12").Value = 47
$ws.Range("Q12").Value = 7
$ws.Range("R12").Value = 7
$ws.Range("S12").Value = 7
$ws.Range("T12").Value = 7
$ws.Range("U12").Value = 7
$ws.Range("V12").Value = 7

# Row 13
$ws.Range("A13").Value = 53
$ws.Range("B13").Value = 'memory'
$ws.Range("C13").Value = 6
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 12
$ws.Range("F13").Value = 376
$ws.Range("G13").Value = 'living_rooms'
$ws.Range("H13").Value = 'living_rooms'
$ws.Range("I13").ClearContents()
$ws.Range("J13").Value = 'new'
$ws.Range("K13").Value = 'f'
$ws.Range("L13").Value = 'stimuli/img_wbws6.png'
$ws.Range("M13").Value = 57.97777777777777
$ws.Range("N13").Value = 42.53333333333333
$ws.Range("O13").Value = 50.25555555555555
$ws.Range("P13").Value = 45
$ws.Range("Q13").Value = 4
$ws.Range("R13").Value = 4
$ws.Range("S13").Value = 4
$ws.Range("T13").Value = 4
$ws.Range("U13").Value = 4
$ws.Range("V13").Value = 5

# Row 14
$ws.Range("A14").Value = 53
$ws.Range("B14").Value = 'memory'
$ws.Range("C14").Value = 6
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 377
$ws.Range("G14").Value = 'living_rooms'
$ws.Range("H14").Value = 'living_rooms'
$ws.Range("I14").Value = 'target'
$ws.Range("J14").Value = 'old'
$ws.Range("K14").Value = 'j'
$ws.Range("L14").Value = 'stimuli/img_dg5h7.png'
$ws.Range("M14").Value = 88.72093023255815
$ws.Range("N14").Value = 76.06976744186046
$ws.Range("O14").Value = 82.3953488372093
$ws.Range("P14").Value = 43
$ws.Range("Q14").Value = 10
$ws.Range("R14").Value = 10
$ws.Range("S14").Value = 10
$ws.Range("T14").Value = 10
$ws.Range("U14").Value = 10
$ws.Range("V14").Value = 10

# Row 15
$ws.Range("A15").Value = 53
$ws.Range("B15").Value = 'memory'
$ws.Range("C15").Value = 6
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 14
$ws.Range("F15").Value = 378
$ws.Range("G15").Value = 'living_rooms'
$ws.Range("H15").Value = 'living_rooms'
$ws.Range("I15").ClearContents()
$ws.Range("J15").Value = 'new'
$ws.Range("K15").Value = 'f'
$ws.Range("L15").Value = 'stimuli/img_q9lab.png'
$ws.Range("M15").Value = 53.97560975609756
$ws.Range("N15").Value = 32.90243902439025
$ws.Range("O15").Value = 43.4390243902439
$ws.Range("P15").Value = 41
$ws.Range("Q15").Value = 3
$ws.Range("R15").Value = 3
$ws.Range("S15").Value = 3
$ws.Range("T15").Value = 3
$ws.Range("U15").Value = 4
$ws.Range("V15").Value = 3

# Row 16
$ws.Range("A16").Value = 53
$ws.Range("B16").Value = 'memory'
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 15
$ws.Range("F16").Value = 379
$ws.Range("G16").Value = 'living_rooms'
$ws.Range("H16").Value = 'living_rooms'
$ws.Range("I16").ClearContents()
$ws.Range("J16").Value = 'new'
$ws.Range("K16").Value = 'f'
$ws.Range("L16").Value = 'stimuli/img_pjfx6.png'
$ws.Range("M16").Value = 32.23404255319149
$ws.Range("N16").Value = 26.59574468085106
$ws.Range("O16").Value = 29.41489361702127
$ws.Range("P16").Value = 47
$ws.Range("Q16").Value = 2
$ws.Range("R16").Value = 2
$ws.Range("S16").Value = 2
$ws.Range("T16").Value = 2
$ws.Range("U16").Value = 2
$ws.Range("V16").Value = 3

# Row 17
$ws.Range("A17").Value = 53
$ws.Range("B17").Value = 'memory'
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = 380
$ws.Range("G17").Value = 'living_rooms'
$ws.Range("H17").Value = 'living_rooms'
$ws.Range("I17").Value = 'target'
$ws.Range("J17").Value = 'old'
$ws.Range("K17").Value = 'j'
$ws.Range("L17").Value = 'stimuli/img_w8yhd.png'
$ws.Range("M17").Value = 55.74418604651163
$ws.Range("N17").Value = 38.90697674418605
$ws.Range("O17").Value = 47.32558139534883
$ws.Range("P17").Value = 43
$ws.Range("Q17").Value = 4
$ws.Range("R17").Value = 4
$ws.Range("S17").Value = 4
$ws.Range("T17").Value = 4
$ws.Range("U17").Value = 4
$ws.Range("V17").Value = 4

# Row 18
$ws.Range("A18").Value = 53
$ws.Range("B18").Value = 'memory'
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 17
$ws.Range("F18").Value = 381
$ws.Range("G18").Value = 'living_rooms'
$ws.Range("H18").Value = 'living_rooms'
$ws.Range("I18").Value = 'target'
$ws.Range("J18").Value = 'old'
$ws.Range("K18").Value = 'j'
$ws.Range("L18").Value = 'stimuli/img_rych7.png'
$ws.Range("M18").Value = 30.4468085106383
$ws.Range("N18").Value = 23.4468085106383
$ws.Range("O18").Value = 26.9468085106383
$ws.Range("P18").Value = 47
$ws.Range("Q18").Value = 2
$ws.Range("R18").Value = 2
$ws.Range("S18").Value = 2
$ws.Range("T18").Value = 2
$ws.Range("U18").Value = 2
$ws.Range("V18").Value = 2

# Row 19
$ws.Range("A19").Value = 53
$ws.Range("B19").Value = 'memory'
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 18
$ws.Range("F19").Value = 382
$ws.Range("G19").Value = 'living_rooms'
$ws.Range("H19").Value = 'living_rooms'
$ws.Range("I19").Value = 'target'
$ws.Range("J19").Value = 'old'
$ws.Range("K19").Value = 'j'
$ws.Range("L19").Value = 'stimuli/img_nb8p4.png'
$ws.Range("M19").Value = 16.36170212765957
$ws.Range("N19").Value = 12.70212765957447
$ws.Range("O19").Value = 14.53191489361702
$ws.Range("P19").Value = 47
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = 1
$ws.Range("S19").Value = 1
$ws.Range("T19").Value = 1
$ws.Range("U19").Value = 1
$ws.Range("V19").Value = 1

# Row 20
$ws.Range("A20").Value = 53
$ws.Range("B20").Value = 'memory'
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 19
$ws.Range("F20").Value = 383
$ws.Range("G20").Value = 'living_rooms'
$ws.Range("H20").Value = 'living_rooms'
$ws.Range("I20").Value = 'target'
$ws.Range("J20").Value = 'old'
$ws.Range("K20").Value = 'j'
$ws.Range("L20").Value = 'stimuli/img_jkm86.png'
$ws.Range("M20").Value = 58.32558139534883
$ws.Range("N20").Value = 38.65116279069768
$ws.Range("O20").Value = 48.48837209302326
$ws.Range("P20").Value = 43
$ws.Range("Q20").Value = 4
$ws.Range("R20").Value = 4
$ws.Range("S20").Value = 4
$ws.Range("T20").Value = 4
$ws.Range("U20").Value = 4
$ws.Range("V20").Value = 4

# Row 21
$ws.Range("A21").Value = 53
$ws.Range("B21").Value = 'memory'
$ws.Range("C21").Value = 6
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 20
$ws.Range("F21").Value = 384
$ws.Range("G21").Value = 'living_rooms'
$ws.Range("H21").Value = 'living_rooms'
$ws.Range("I21").ClearContents()
$ws.Range("J21").Value = 'new'
$ws.Range("K21").Value = 'f'
$ws.Range("L21").Value = 'stimuli/img_vgaye.png'
$ws.Range("M21").Value = 80.33333333333333
$ws.Range("N21").Value = 64.57777777777778
$ws.Range("O21").Value = 72.45555555555555
$ws.Range("P21").Value = 45
$ws.Range("Q21").Value = 8
$ws.Range("R21").Value = 8
$ws.Range("S21").Value = 8
$ws.Range("T21").Value = 8
$ws.Range("U21").Value = 7
$ws.Range("V21").Value = 8

# Row 22
$ws.Range("A22").Value = 53
$ws.Range("B22").Value = 'memory'
$ws.Range("C22").Value = 6
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 21
$ws.Range("F22").Value = 385
$ws.Range("G22").Value = 'living_rooms'
$ws.Range("H22").Value = 'living_rooms'
$ws.Range("I22").ClearContents()
$ws.Range("J22").Value = 'new'
$ws.Range("K22").Value = 'f'
$ws.Range("L22").Value = 'stimuli/img_xr3up.png'
$ws.Range("M22").Value = 76.24444444444444
$ws.Range("N22").Value = 55.88888888888889
$ws.Range("O22").Value = 66.06666666666666
$ws.Range("P22").Value = 45
$ws.Range("Q22").Value = 7
$ws.Range("R22").Value = 7
$ws.Range("S22").Value = 7
$ws.Range("T22").Value = 6
$ws.Range("U22").Value = 6
$ws.Range("V22").Value = 6

# Row 23
$ws.Range("A23").Value = 53
$ws.Range("B23").Value = 'memory'
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 22
$ws.Range("F23").Value = 386
$ws.Range("G23").Value = 'living_rooms'
$ws.Range("H23").ClearContents()
$ws.Range("I23").ClearContents()
$ws.Range("J23").Value = 'catch'
$ws.Range("K23").Value = 'f'
$ws.Range("L23").Value = 'stimuli/catch_21.jpg'
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("O23").ClearContents()
$ws.Range("P23").ClearContents()
$ws.Range("Q23").ClearContents()
$ws.Range("R23").ClearContents()
$ws.Range("S23").ClearContents()
$ws.Range("T23").ClearContents()
$ws.Range("U23").ClearContents()
$ws.Range("V23").ClearContents()

# Row 24
$ws.Range("A24").Value = 53
$ws.Range("B24").Value = 'memory'
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 23
$ws.Range("F24").Value = 387
$ws.Range("G24").Value = 'living_rooms'
$ws.Range("H24").Value = 'living_rooms'
$ws.Range("I24").Value = 'target'
$ws.Range("J24").Value = 'old'
$ws.Range("K24").Value = 'j'
$ws.Range("L24").Value = 'stimuli/img_3sw8t.png'
$ws.Range("M24").Value = 67.4888888888889
$ws.Range("N24").Value = 48.51111111111111
$ws.Range("O24").Value = 58
$ws.Range("P24").Value = 45
$ws.Range("Q24").Value = 5
$ws.Range("R24").Value = 5
$ws.Range("S24").Value = 5
$ws.Range("T24").Value = 5
$ws.Range("U24").Value = 5
$ws.Range("V24").Value = 5

# Row 25
$ws.Range("A25").Value = 53
$ws.Range("B25").Value = 'memory'
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 24
$ws.Range("F25").Value = 388
$ws.Range("G25").Value = 'living_rooms'
$ws.Range("H25").Value = 'living_rooms'
$ws.Range("I25").Value = 'target'
$ws.Range("J25").Value = 'old'
$ws.Range("K25").Value = 'j'
$ws.Range("L25").Value = 'stimuli/img_xzyzy.png'
$ws.Range("M25").Value = 85.37209302325581
$ws.Range("N25").Value = 68.90697674418605
$ws.Range("O25").Value = 77.13953488372093
$ws.Range("P25").Value = 43
$ws.Range("Q25").Value = 9
$ws.Range("R25").Value = 9
$ws.Range("S25").Value = 9
$ws.Range("T25").Value = 9
$ws.Range("U25").Value = 9
$ws.Range("V25").Value = 9

# Row 26
$ws.Range("A26").Value = 53
$ws.Range("B26").Value = 'memory'
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 25
$ws.Range("F26").Value = 389
$ws.Range("G26").Value = 'living_rooms'
$ws.Range("H26").Value = 'living_rooms'
$ws.Range("I26").ClearContents()
$ws.Range("J26").Value = 'new'
$ws.Range("K26").Value = 'f'
$ws.Range("L26").Value = 'stimuli/img_hmmra.png'
$ws.Range("M26").Value = 54.65853658536585
$ws.Range("N26").Value = 34.24390243902439
$ws.Range("O26").Value = 44.45121951219512
$ws.Range("P26").Value = 41
$ws.Range("Q26").Value = 3
$ws.Range("R26").Value = 3
$ws.Range("S26").Value = 3
$ws.Range("T26").Value = 4
$ws.Range("U26").Value = 4
$ws.Range("V26").Value = 3

# Row 27
$ws.Range("A27").Value = 53
$ws.Range("B27").Value = 'memory'
$ws.Range("C27").Value = 6
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 26
$ws.Range("F27").Value = 390
$ws.Range("G27").Value = 'living_rooms'
$ws.Range("H27").Value = 'living_rooms'
$ws.Range("I27").ClearContents()
$ws.Range("J27").Value = 'new'
$ws.Range("K27").Value = 'f'
$ws.Range("L27").Value = 'stimuli/img_tn8ys.png'
$ws.Range("M27").Value = 86.70454545454545
$ws.Range("N27").Value = 72.4090909090909
$ws.Range("O27").Value = 79.55681818181819
$ws.Range("P27").Value = 44
$ws.Range("Q27").Value = 10
$ws.Range("R27").Value = 10
$ws.Range("S27").Value = 10
$ws.Range("T27").Value = 9
$ws.Range("U27").Value = 9
$ws.Range("V27").Value = 10

# Row 28
$ws.Range("A28").Value = 53
$ws.Range("B28").Value = 'memory'
$ws.Range("C28").Value = 6
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 27
$ws.Range("F28").Value = 391
$ws.Range("G28").Value = 'living_rooms'
$ws.Range("H28").Value = 'living_rooms'
$ws.Range("I28").Value = 'target'
$ws.Range("J28").Value = 'old'
$ws.Range("K28").Value = 'j'
$ws.Range("L28").Value = 'stimuli/img_o30wb.png'
$ws.Range("M28").Value = 81.06666666666666
$ws.Range("N28").Value = 65.37777777777778
$ws.Range("O28").Value = 73.22222222222223
$ws.Range("P28").Value = 45
$ws.Range("Q28").Value = 8
$ws.Range("R28").Value = 8
$ws.Range("S28").Value = 8
$ws.Range("T28").Value = 8
$ws.Range("U28").Value = 8
$ws.Range("V28").Value = 8

# Row 29
$ws.Range("A29").Value = 53
$ws.Range("B29").Value = 'memory'
$ws.Range("C29").Value = 6
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 28
$ws.Range("F29").Value = 392
$ws.Range("G29").Value = 'living_rooms'
$ws.Range("H29").Value = 'living_rooms'
$ws.Range("I29").ClearContents()
$ws.Range("J29").Value = 'new'
$ws.Range("K29").Value = 'f'
$ws.Range("L29").Value = 'stimuli/img_pdzf1.png'
$ws.Range("M29").Value = 86.23913043478261
$ws.Range("N29").Value = 67.17391304347827
$ws.Range("O29").Value = 76.70652173913044
$ws.Range("P29").Value = 46
$ws.Range("Q29").Value = 9
$ws.Range("R29").Value = 9
$ws.Range("S29").Value = 9
$ws.Range("T29").Value = 9
$ws.Range("U29").Value = 9
$ws.Range("V29").Value = 8

# Row 30
$ws.Range("A30").Value = 53
$ws.Range("B30").Value = 'memory'
$ws.Range("C30").Value = 6
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 29
$ws.Range("F30").Value = 393
$ws.Range("G30").Value = 'living_rooms'
$ws.Range("H30").Value = 'living_rooms'
$ws.Range("I30").Value = 'target'
$ws.Range("J30").Value = 'old'
$ws.Range("K30").Value = 'j'
$ws.Range("L30").Value = 'stimuli/img_8dmpq.png'
$ws.Range("M30").Value = 30.65909090909091
$ws.Range("N30").Value = 24.11363636363636
$ws.Range("O30").Value = 27.38636363636364
$ws.Range("P30").Value = 44
$ws.Range("Q30").Value = 2
$ws.Range("R30").Value = 2
$ws.Range("S30").Value = 2
$ws.Range("T30").Value = 2
$ws.Range("U30").Value = 2
$ws.Range("V30").Value = 2

# Row 31
$ws.Range("A31").Value = 53
$ws.Range("B31").Value = 'memory'
$ws.Range("C31").Value = 6
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = 30
$ws.Range("F31").Value = 394
$ws.Range("G31").Value = 'living_rooms'
$ws.Range("H31").Value = 'living_rooms'
$ws.Range("I31").ClearContents()
$ws.Range("J31").Value = 'new'
$ws.Range("K31").Value = 'f'
$ws.Range("L31").Value = 'stimuli/img_gka64.png'
$ws.Range("M31").Value = 19.23809523809524
$ws.Range("N31").Value = 20.02380952380953
$ws.Range("O31").Value = 19.63095238095238
$ws.Range("P31").Value = 42
$ws.Range("Q31").Value = 1
$ws.Range("R31").Value = 1
$ws.Range("S31").Value = 1
$ws.Range("T31").Value = 1
$ws.Range("U31").Value = 1
$ws.Range("V31").Value = 2

# Row 32
$ws.Range("A32").Value = 53
$ws.Range("B32").Value = 'memory'
$ws.Range("C32").Value = 6
$ws.Range("D32").Value = 1
$ws.Range("E32").Value = 31
$ws.Range("F32").Value = 395
$ws.Range("G32").Value = 'living_rooms'
$ws.Range("H32").Value = 'living_rooms'
$ws.Range("I32").ClearContents()
$ws.Range("J32").Value = 'new'
$ws.Range("K32").Value = 'f'
$ws.Range("L32").Value = 'stimuli/img_ra2nm.png'
$ws.Range("M32").Value = 70.75
$ws.Range("N32").Value = 50.375
$ws.Range("O32").Value = 60.5625
$ws.Range("P32").Value = 40
$ws.Range("Q32").Value = 6
$ws.Range("R32").Value = 6
$ws.Range("S32").Value = 6
$ws.Range("T32").Value = 5
$ws.Range("U32").Value = 5
$ws.Range("V32").Value = 6

# Row 33
$ws.Range("A33").Value = 53
$ws.Range("B33").Value = 'memory'
$ws.Range("C33").Value = 6
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = 32
$ws.Range("F33").Value = 396
$ws.Range("G33").Value = 'living_rooms'
$ws.Range("H33").Value = 'living_rooms'
$ws.Range("I33").Value = 'target'
$ws.Range("J33").Value = 'old'
$ws.Range("K33").Value = 'j'
$ws.Range("L33").Value = 'stimuli/img_g13d5.png'
$ws.Range("M33").Value = 73
$ws.Range("N33").Value = 51.51111111111111
$ws.Range("O33").Value = 62.25555555555556
$ws.Range("P33").Value = 45
$ws.Range("Q33").Value = 6
$ws.Range("R33").Value = 6
$ws.Range("S33").Value = 6
$ws.Range("T33").Value = 6
$ws.Range("U33").Value = 6
$ws.Range("V33").Value = 6

# Row 34
$ws.Range("A34").Value = 53
$ws.Range("B34").Value = 'memory'
$ws.Range("C34").Value = 6
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = 33
$ws.Range("F34").Value = 397
$ws.Range("G34").Value = 'living_rooms'
$ws.Range("H34").Value = 'living_rooms'
$ws.Range("I34").Value = 'target'
$ws.Range("J34").Value = 'old'
$ws.Range("K34").Value = 'j'
$ws.Range("L34").Value = 'stimuli/img_c0vzo.png'
$ws.Range("M34").Value = 21.51162790697675
$ws.Range("N34").Value = 8.232558139534884
$ws.Range("O34").Value = 14.87209302325581
$ws.Range("P34").Value = 43
$ws.Range("Q34").Value = 1
$ws.Range("R34").Value = 1
$ws.Range("S34").Value = 1
$ws.Range("T34").Value = 1
$ws.Range("U34").Value = 1
$ws.Range("V34").Value = 1

# Row 35
$ws.Range("A35").Value = 53
$ws.Range("B35").Value = 'memory'
$ws.Range("C35").Value = 6
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = 34
$ws.Range("F35").Value = 398
$ws.Range("G35").Value = 'living_rooms'
$ws.Range("H35").Value = 'living_rooms'
$ws.Range("I35").ClearContents()
$ws.Range("J35").Value = 'new'
$ws.Range("K35").Value = 'f'
$ws.Range("L35").Value = 'stimuli/img_z4jxm.png'
$ws.Range("M35").Value = 88.30952380952381
$ws.Range("N35").Value = 72.64285714285714
$ws.Range("O35").Value = 80.47619047619048
$ws.Range("P35").Value = 42
$ws.Range("Q35").Value = 10
$ws.Range("R35").Value = 10
$ws.Range("S35").Value = 10
$ws.Range("T35").Value = 10
$ws.Range("U35").Value = 10
$ws.Range("V35").Value = 10

# Row 36
$ws.Range("A36").Value = 53
$ws.Range("B36").Value = 'memory'
$ws.Range("C36").Value = 6
$ws.Range("D36").Value = 1
$ws.Range("E36").Value = 35
$ws.Range("F36").Value = 399
$ws.Range("G36").Value = 'living_rooms'
$ws.Range("H36").Value = 'living_rooms'
$ws.Range("I36").Value = 'target'
$ws.Range("J36").Value = 'old'
$ws.Range("K36").Value = 'j'
$ws.Range("L36").Value = 'stimuli/img_zxvl3.png'
$ws.Range("M36").Value = 68.78260869565217
$ws.Range("N36").Value = 47.56521739130435
$ws.Range("O36").Value = 58.17391304347827
$ws.Range("P36").Value = 46
$ws.Range("Q36").Value = 5
$ws.Range("R36").Value = 5
$ws.Range("S36").Value = 5
$ws.Range("T36").Value = 5
$ws.Range("U36").Value = 5
$ws.Range("V36").Value = 5

# Row 37
$ws.Range("A37").Value = 53
$ws.Range("B37").Value = 'memory'
$ws.Range("C37").Value = 6
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 36
$ws.Range("F37").Value = 400
$ws.Range("G37").Value = 'living_rooms'
$ws.Range("H37").Value = 'living_rooms'
$ws.Range("I37").Value = 'target'
$ws.Range("J37").Value = 'old'
$ws.Range("K37").Value = 'j'
$ws.Range("L37").Value = 'stimuli/img_jpldg.png'
$ws.Range("M37").Value = 79.54545454545455
$ws.Range("N37").Value = 57.75
$ws.Range("O37").Value = 68.64772727272728
$ws.Range("P37").Value = 44
$ws.Range("Q37").Value = 7
$ws.Range("R37").Value = 7
$ws.Range("S37").Value = 7
$ws.Range("T37").Value = 7
$ws.Range("U37").Value = 7
$ws.Range("V37").Value = 7

# Row 38
$ws.Range("A38").Value = 53
$ws.Range("B38").Value = 'memory'
$ws.Range("C38").Value = 6
$ws.Range("D38").Value = 1
$ws.Range("E38").Value = 37
$ws.Range("F38").Value = 401
$ws.Range("G38").Value = 'living_rooms'
$ws.Range("H38").Value = 'living_rooms'
$ws.Range("I38").Value = 'target'
$ws.Range("J38").Value = 'old'
$ws.Range("K38").Value = 'j'
$ws.Range("L38").Value = 'stimuli/img_165pk.png'
$ws.Range("M38").Value = 85.73333333333333
$ws.Range("N38").Value = 69.22222222222223
$ws.Range("O38").Value = 77.47777777777779
$ws.Range("P38").Value = 45
$ws.Range("Q38").Value = 9
$ws.Range("R38").Value = 9
$ws.Range("S38").Value = 9
$ws.Range("T38").Value = 9
$ws.Range("U38").Value = 9
$ws.Range("V38").Value = 9

# Row 39
$ws.Range("A39").Value = 53
$ws.Range("B39").Value = 'memory'
$ws.Range("C39").Value = 6
$ws.Range("D39").Value = 1
$ws.Range("E39").Value = 38
$ws.Range("F39").Value = 402
$ws.Range("G39").Value = 'living_rooms'
$ws.Range("H39").Value = 'living_rooms'
$ws.Range("I39").Value = 'target'
$ws.Range("J39").Value = 'old'
$ws.Range("K39").Value = 'j'
$ws.Range("L39").Value = 'stimuli/img_5jy9c.png'
$ws.Range("M39").Value = 87.37209302325581
$ws.Range("N39").Value = 79.18604651162791
$ws.Range("O39").Value = 83.27906976744185
$ws.Range("P39").Value = 43
$ws.Range("Q39").Value = 10
$ws.Range("R39").Value = 10
$ws.Range("S39").Value = 10
$ws.Range("T39").Value = 10
$ws.Range("U39").Value = 9
$ws.Range("V39").Value = 10

# Row 40
$ws.Range("A40").Value = 53
$ws.Range("B40").Value = 'memory'
$ws.Range("C40").Value = 6
$ws.Range("D40").Value = 1
$ws.Range("E40").Value = 39
$ws.Range("F40").Value = 403
$ws.Range("G40").Value = 'living_rooms'
$ws.Range("H40").Value = 'living_rooms'
$ws.Range("I40").Value = 'target'
$ws.Range("J40").Value = 'old'
$ws.Range("K40").Value = 'j'
$ws.Range("L40").Value = 'stimuli/img_9bkl9.png'
$ws.Range("M40").Value = 46.62162162162162
$ws.Range("N40").Value = 34.27027027027027
$ws.Range("O40").Value = 40.44594594594595
$ws.Range("P40").Value = 37
$ws.Range("Q40").Value = 3
$ws.Range("R40").Value = 3
$ws.Range("S40").Value = 3
$ws.Range("T40").Value = 3
$ws.Range("U40").Value = 3
$ws.Range("V40").Value = 3

# Row 41
$ws.Range("A41").Value = 53
$ws.Range("B41").Value = 'memory'
$ws.Range("C41").Value = 6
$ws.Range("D41").Value = 1
$ws.Range("E41").Value = 40
$ws.Range("F41").Value = 404
$ws.Range("G41").Value = 'living_rooms'
$ws.Range("H41").Value = 'living_rooms'
$ws.Range("I41").ClearContents()
$ws.Range("J41").Value = 'new'
$ws.Range("K41").Value = 'f'
$ws.Range("L41").Value = 'stimuli/img_b21d7.png'
$ws.Range("M41").Value = 27.75555555555556
$ws.Range("N41").Value = 13.86666666666667
$ws.Range("O41").Value = 20.81111111111111
$ws.Range("P41").Value = 45
$ws.Range("Q41").Value = 1
$ws.Range("R41").Value = 1
$ws.Range("S41").Value = 1
$ws.Range("T41").Value = 2
$ws.Range("U41").Value = 2
$ws.Range("V41").Value = 1

# Row 42
$ws.Range("A42").Value = 53
$ws.Range("B42").Value = 'memory'
$ws.Range("C42").Value = 6
$ws.Range("D42").Value = 1
$ws.Range("E42").Value = 41
$ws.Range("F42").Value = 405
$ws.Range("G42").Value = 'living_rooms'
$ws.Range("H42").Value = 'living_rooms'
$ws.Range("I42").ClearContents()
$ws.Range("J42").Value = 'new'
$ws.Range("K42").Value = 'f'
$ws.Range("L42").Value = 'stimuli/img_16kib.png'
$ws.Range("M42").Value = 80.97727272727273
$ws.Range("N42").Value = 61.11363636363637
$ws.Range("O42").Value = 71.04545454545455
$ws.Range("P42").Value = 44
$ws.Range("Q42").Value = 8
$ws.Range("R42").Value = 8
$ws.Range("S42").Value = 8
$ws.Range("T42").Value = 7
$ws.Range("U42").Value = 7
$ws.Range("V42").Value = 7

